$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new columns before column BC (shifts old BC.. right by 2)
$ws.Range("BC1:BD1").EntireColumn.Insert()

# Match column width of the new columns to their left neighbour (BB)
$bbWidth = $ws.Range("BB1").ColumnWidth
$ws.Range("BC1").ColumnWidth = $bbWidth
$ws.Range("BD1").ColumnWidth = $bbWidth

# Copy header formatting from an existing Name/ID-style merged pair (AA4:AB4)
# onto the new BC4:BD4 pair, then merge + fill in the new "Crumb" header
$ws.Range("AA4:AB4").Copy()
$ws.Range("BC4:BD4").PasteSpecial(-4122)
$ws.Range("BC4:BD4").Merge()

$ws.Range("BC4").Value = "Crumb"
$ws.Range("BC5").Value = "Name"
$ws.Range("BD5").Value = "ID"
$ws.Range("BC6").Value = "Crumb"
$ws.Range("BD6").Value = 15000

# Restore selection to match the saved view state
$ws.Range("BE9").Select()
